$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Category" header in column H
$ws.Range("H1").Value = "Category"

# Fill in the previously-empty ISIN / Amount / First trading data for 4iG Nyrt. (row 2)
$ws.Range("D2").Value = "HU0000167788"
$ws.Range("E2").Value = 94000000
$ws.Range("F2").Value = "2004.09.22."
$ws.Range("D2:F2").WrapText = $true

# Populate the Category column.
# Rows 13-21 (Budapesti Elektromos Muvek .. SET GROUP) are "B" shares,
# everything else is "A" shares.
for ($r = 2; $r -le 30; $r++) {
    if ($r -ge 13 -and $r -le 21) {
        $ws.Cells.Item($r, 8).Value = "W_RESZVENYB"
    } else {
        $ws.Cells.Item($r, 8).Value = "W_RESZVENYA"
    }
}

$ws.Range("J9").Select()
